$wb = $excel.ActiveWorkbook

# Reorder worksheets: "review_info" should come before "hotel_info"
$hotelInfo = $wb.Worksheets.Item("hotel_info")
$reviewInfo = $wb.Worksheets.Item("review_info")
$reviewInfo.Move($hotelInfo)

# Insert a new "State" column into hotel_info, between Hotel_Name and City
$ws = $wb.Worksheets.Item("hotel_info")
$ws.Columns("C").Insert()
$ws.Range("C1").Value = "State"
$ws.Range("C2").Value = "Louisiana"
